$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A75").Value = "30-11-2025"
$ws.Range("B75").Value = "The price of gold in India today is ₹12,982 per gram for 24 karat gold, ₹11,900 per gram for 22 karat gold and ₹9,737 per gram for 18 karat gold (also called 999 gold)."
